$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "arrray" -> "array" (row 30, column B)
$ws.Range("B30").Value = "array"

# New row 42: search in rotated sorted array
$ws.Range("A42").Value = "search in rotated sorted array"
$ws.Range("B42").Value = "array"
$ws.Range("C42").Value = "two pointers"
$ws.Range("E42").Value = 45
$ws.Range("F42").Value = "https://leetcode.com/problems/search-in-rotated-sorted-array/"

# New row 43: combination sum (hyperlink/url entered before the rest of the row)
$ws.Range("F43").Value = "https://leetcode.com/problems/combination-sum/"
$ws.Range("A43").Value = "combination sum"
$ws.Range("B43").Value = "array"
$ws.Range("C43").Value = "dfs"
$ws.Range("D43").Value = "backtracking"
$ws.Range("E43").Value = 25

# New row 44: permutations (hyperlink/url entered before the rest of the row)
$ws.Range("F44").Value = "https://leetcode.com/problems/permutations/"
$ws.Range("A44").Value = "permutations"
$ws.Range("B44").Value = "array"
$ws.Range("C44").Value = "recursive"
$ws.Range("D44").Value = "dfs"
$ws.Range("E44").Value = 31

# New row 45: merge intervals (hyperlink/url entered before the rest of the row)
$ws.Range("F45").Value = "https://leetcode.com/problems/merge-intervals/"
$ws.Range("A45").Value = "merge intervals"
$ws.Range("B45").Value = "sorting"
$ws.Range("C45").Value = "intervals"
$ws.Range("E45").Value = 12

# New row 46: lowest common ancestor of a binary tree
$ws.Range("A46").Value = "lowest common ancestor of a binary tree"
$ws.Range("F46").Value = "https://leetcode.com/problems/lowest-common-ancestor-of-a-binary-tree/"

# Match the existing "url" column formatting used elsewhere on the sheet
# (F31..F35/F37/F39 carry style index 1 - Consolas/vertical-centered look);
# copy that formatting onto the new F42:F44 cells instead of redefining fonts/styles.
$ws.Range("F31").Copy()
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("F43").PasteSpecial(-4122)
$ws.Range("F44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update view state to match target
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("A47").Select()
